# Add the Eadie-Hofstee style auxiliary table (columns D:F) to "Hoja1",
# used to compute enzyme kinematic parameters (Lineweaver-Burk / Eadie-Hofstee
# / Hanes-Woolf representations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row headers for the new "a"/"b" parameters
$ws.Range("D2").Value = "a"
$ws.Range("D3").Value = "b"

# Computed parameter values
$ws.Range("E2").Value = 0.0594814
$ws.Range("E3").Value = -0.0279608

# "2*" column header and formulas doubling the parameter values
$ws.Range("F1").Value = "2*"
$ws.Range("F2").Formula = "=2*E2"
$ws.Range("F3").Formula = "=2*E3"

# Leave the selection on the last entered formula cell
$ws.Range("F3").Select() | Out-Null
